$d = $word.ActiveDocument

# Step 0: update the "Curso (semestre ideal)" line
$d.Content.Find.Execute("Curso (semestre ideal): EQN (12)", $true, $false, $false, $false, $false, $true, 1, $false, "Curso (semestre ideal): EQD (10), EQN (12)", 2) | Out-Null

# Step 1: replace each Requisitos line with a unique placeholder token (avoids ambiguous matches
# once values start repeating across old/new positions during the reorder).
$d.Content.Find.Execute("LOB1003 -  Cálculo I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT00@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1004 -  Cálculo II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT01@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1006 -  Cálculo IV  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT02@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT03@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1011 -  Eletricidade Aplicada  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT04@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1012 -  Estatística  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT05@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1018 -  Física I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT06@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1019 -  Física II  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT07@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1024 -  Mecânica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT08@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1036 -  Geometria Analítica  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT09@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1037 -  Àlgebra Linear  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT10@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1038 -  Física Experimental I  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT11@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1039 -  Física Experimental III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT12@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1040 -  Laboratório de Eletricidade  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT13@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT14@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1052 -  Cálculo III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT15@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1053 -  Física III  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT16@@", 2) | Out-Null
$d.Content.Find.Execute("LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT17@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT18@@", 2) | Out-Null
$d.Content.Find.Execute("LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT19@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT20@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4095 -  Química Geral Experimental  (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT21@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT22@@", 2) | Out-Null
$d.Content.Find.Execute("LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)", $true, $false, $false, $false, $false, $true, 1, $false, "@@SLOT23@@", 2) | Out-Null

# Step 2: replace each placeholder token with the final text for that slot.
$d.Content.Find.Execute("@@SLOT00@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT01@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1004 -  Cálculo II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT02@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1006 -  Cálculo IV  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT03@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT04@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1011 -  Eletricidade Aplicada  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT05@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1012 -  Estatística  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT06@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1018 -  Física I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT07@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT08@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4010 -  Introdução à  Engenharia  Química  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT09@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1039 -  Física Experimental III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT10@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1040 -  Laboratório de Eletricidade  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT11@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1052 -  Cálculo III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT12@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1053 -  Física III  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT13@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT14@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT15@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4095 -  Química Geral Experimental  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT16@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1019 -  Física II  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT17@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1024 -  Mecânica  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT18@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1036 -  Geometria Analítica  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT19@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1038 -  Física Experimental I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT20@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1037 -  Àlgebra Linear  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT21@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT22@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOB1003 -  Cálculo I  (Requisito)", 2) | Out-Null
$d.Content.Find.Execute("@@SLOT23@@", $true, $false, $false, $false, $false, $true, 1, $false, "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)", 2) | Out-Null

Write-Host "done"